# Auto-generated edit script applying market-price / profit-column updates
# to the Shinryu_Profits workbook (8 Leves sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 3310.15
$ws.Cells.Item(74, 9).Value = 3301.5
$ws.Cells.Item(74, 10).Value = 3311.111
$ws.Cells.Item(74, 11).Value = 3301.5
$ws.Cells.Item(74, 12).Value = 3311.111
$ws.Cells.Item(74, 13).Value = -2365.5
$ws.Cells.Item(74, 14).Value = -5183.111
$ws.Cells.Item(77, 8).Value = 3310.15
$ws.Cells.Item(77, 9).Value = 3301.5
$ws.Cells.Item(77, 10).Value = 3311.111
$ws.Cells.Item(77, 11).Value = 16507.5
$ws.Cells.Item(77, 12).Value = 16555.555
$ws.Cells.Item(77, 13).Value = -11827.5
$ws.Cells.Item(77, 14).Value = -25915.555
$ws.Cells.Item(87, 8).Value = 32119
$ws.Cells.Item(87, 10).Value = 37542.8
$ws.Cells.Item(87, 12).Value = 37542.8
$ws.Cells.Item(87, 14).Value = -40038.8
$ws.Cells.Item(90, 8).Value = 32119
$ws.Cells.Item(90, 10).Value = 37542.8
$ws.Cells.Item(90, 12).Value = 112628.4
$ws.Cells.Item(90, 14).Value = -125108.4
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 14).ClearContents()  # N97 removed entirely in target state
$ws.Cells.Item(112, 8).Value = 1271.6052
$ws.Cells.Item(112, 10).Value = 1323.4857
$ws.Cells.Item(112, 12).Value = 3970.4571
$ws.Cells.Item(112, 14).Value = -6186.4571
$ws.Cells.Item(129, 8).Value = 1100.46
$ws.Cells.Item(129, 9).Value = 565.6667
$ws.Cells.Item(129, 11).Value = 1697.0001
$ws.Cells.Item(129, 13).Value = 3302.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 19234598
$ws.Cells.Item(32, 9).Value = 22730444
$ws.Cells.Item(32, 11).Value = 22730444
$ws.Cells.Item(32, 13).Value = -22730157
$ws.Cells.Item(37, 8).Value = 12248.25
$ws.Cells.Item(37, 9).Value = 3496.5
$ws.Cells.Item(37, 10).Value = 21000
$ws.Cells.Item(37, 11).Value = 3496.5
$ws.Cells.Item(37, 12).Value = 21000
$ws.Cells.Item(37, 13).Value = -3223.5
$ws.Cells.Item(37, 14).Value = -21546
$ws.Cells.Item(44, 8).Value = 23024.5
$ws.Cells.Item(44, 10).Value = 23024.5
$ws.Cells.Item(44, 12).Value = 23024.5
$ws.Cells.Item(44, 14).Value = -24000.5
$ws.Cells.Item(45, 8).Value = 1565.8
$ws.Cells.Item(45, 9).Value = 867.5
$ws.Cells.Item(45, 10).Value = 2031.3334
$ws.Cells.Item(45, 11).Value = 867.5
$ws.Cells.Item(45, 12).Value = 2031.3334
$ws.Cells.Item(45, 13).Value = -490.5
$ws.Cells.Item(45, 14).Value = -2785.3334
$ws.Cells.Item(74, 8).Value = 1258.091
$ws.Cells.Item(74, 9).Value = 1243.1765
$ws.Cells.Item(74, 10).Value = 1308.8
$ws.Cells.Item(74, 11).Value = 1243.1765
$ws.Cells.Item(74, 12).Value = 1308.8
$ws.Cells.Item(74, 13).Value = -369.1765
$ws.Cells.Item(74, 14).Value = -3056.8
$ws.Cells.Item(77, 8).Value = 1258.091
$ws.Cells.Item(77, 9).Value = 1243.1765
$ws.Cells.Item(77, 10).Value = 1308.8
$ws.Cells.Item(77, 11).Value = 6215.8825
$ws.Cells.Item(77, 12).Value = 6544
$ws.Cells.Item(77, 13).Value = -1847.8825
$ws.Cells.Item(77, 14).Value = -15280
$ws.Cells.Item(80, 8).Value = 25400
$ws.Cells.Item(80, 10).Value = 25400
$ws.Cells.Item(80, 12).Value = 25400
$ws.Cells.Item(80, 14).Value = -27396
$ws.Cells.Item(83, 8).Value = 25400
$ws.Cells.Item(83, 10).Value = 25400
$ws.Cells.Item(83, 12).Value = 76200
$ws.Cells.Item(83, 14).Value = -86184
$ws.Cells.Item(112, 8).Value = 23177.4
$ws.Cells.Item(112, 10).Value = 23177.4
$ws.Cells.Item(112, 12).Value = 23177.4
$ws.Cells.Item(112, 14).Value = -26131.4
$ws.Cells.Item(114, 8).Value = 19933
$ws.Cells.Item(114, 10).Value = 19933
$ws.Cells.Item(114, 12).Value = 19933
$ws.Cells.Item(114, 14).Value = -28611
$ws.Cells.Item(124, 8).Value = 3733.6667
$ws.Cells.Item(124, 10).Value = 3733.6667
$ws.Cells.Item(124, 12).Value = 3733.6667
$ws.Cells.Item(124, 14).Value = -13553.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 21270.908
$ws.Cells.Item(82, 9).Value = 12160
$ws.Cells.Item(82, 10).Value = 24687.5
$ws.Cells.Item(82, 11).Value = 12160
$ws.Cells.Item(82, 12).Value = 24687.5
$ws.Cells.Item(82, 13).Value = -11777
$ws.Cells.Item(82, 14).Value = -25453.5
$ws.Cells.Item(85, 8).Value = 21270.908
$ws.Cells.Item(85, 9).Value = 12160
$ws.Cells.Item(85, 10).Value = 24687.5
$ws.Cells.Item(85, 11).Value = 12160
$ws.Cells.Item(85, 12).Value = 24687.5
$ws.Cells.Item(85, 13).Value = -10834
$ws.Cells.Item(85, 14).Value = -27339.5
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 14).ClearContents()  # N112 removed entirely in target state
$ws.Cells.Item(134, 8).Value = 1519.6604
$ws.Cells.Item(134, 9).Value = 1452.4255
$ws.Cells.Item(134, 10).Value = 2046.3334
$ws.Cells.Item(134, 11).Value = 4357.2765
$ws.Cells.Item(134, 12).Value = 6139.0002
$ws.Cells.Item(134, 13).Value = -1822.2765
$ws.Cells.Item(134, 14).Value = -11209.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4756.857
$ws.Cells.Item(16, 9).Value = 3433.3333
$ws.Cells.Item(16, 10).Value = 5749.5
$ws.Cells.Item(16, 11).Value = 3433.3333
$ws.Cells.Item(16, 12).Value = 5749.5
$ws.Cells.Item(16, 13).Value = -3146.3333
$ws.Cells.Item(16, 14).Value = -6323.5
$ws.Cells.Item(31, 8).Value = 1782
$ws.Cells.Item(31, 9).Value = 1374.762
$ws.Cells.Item(31, 11).Value = 1374.762
$ws.Cells.Item(31, 13).Value = -1079.762
$ws.Cells.Item(34, 8).Value = 1782
$ws.Cells.Item(34, 9).Value = 1374.762
$ws.Cells.Item(34, 10).Value = 3682.4443
$ws.Cells.Item(34, 11).Value = 1374.762
$ws.Cells.Item(34, 13).Value = -1172.762
$ws.Cells.Item(113, 8).Value = 4756.857
$ws.Cells.Item(113, 9).Value = 3433.3333
$ws.Cells.Item(113, 10).Value = 5749.5
$ws.Cells.Item(113, 11).Value = 3433.3333
$ws.Cells.Item(113, 12).Value = 5749.5
$ws.Cells.Item(113, 13).Value = -1263.3333
$ws.Cells.Item(113, 14).Value = -10089.5
$ws.Cells.Item(132, 8).Value = 1794.125
$ws.Cells.Item(132, 9).Value = 1370.9286
$ws.Cells.Item(132, 11).Value = 4112.7858
$ws.Cells.Item(132, 13).Value = -1582.7858
$ws.Cells.Item(134, 8).Value = 1638.7106
$ws.Cells.Item(134, 9).Value = 908.34375
$ws.Cells.Item(134, 10).Value = 5534
$ws.Cells.Item(134, 11).Value = 2725.03125
$ws.Cells.Item(134, 12).Value = 16602
$ws.Cells.Item(134, 13).Value = -190.03125
$ws.Cells.Item(134, 14).Value = -21672

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(97, 8).Value = 679.7143
$ws.Cells.Item(97, 9).Value = 500
$ws.Cells.Item(97, 10).Value = 751.6
$ws.Cells.Item(97, 11).Value = 1500
$ws.Cells.Item(97, 12).Value = 2254.8
$ws.Cells.Item(97, 13).Value = -1004
$ws.Cells.Item(97, 14).Value = -3246.8
$ws.Cells.Item(98, 8).Value = 715.7857
$ws.Cells.Item(98, 9).Value = 594.7143
$ws.Cells.Item(98, 10).Value = 836.8570999999999
$ws.Cells.Item(98, 11).Value = 1784.1429
$ws.Cells.Item(98, 12).Value = 2510.5713
$ws.Cells.Item(98, 13).Value = -286.1428999999998
$ws.Cells.Item(98, 14).Value = -5506.5713
$ws.Cells.Item(122, 8).Value = 8363.786
$ws.Cells.Item(122, 9).Value = 10272.137
$ws.Cells.Item(122, 11).Value = 92449.23300000001
$ws.Cells.Item(122, 13).Value = -89999.23300000001
$ws.Cells.Item(132, 8).Value = 723070.2
$ws.Cells.Item(132, 9).Value = 1094.8
$ws.Cells.Item(132, 11).Value = 9853.199999999999
$ws.Cells.Item(132, 13).Value = -7323.199999999999
$ws.Cells.Item(137, 8).Value = 2966.353
$ws.Cells.Item(137, 10).Value = 4100
$ws.Cells.Item(137, 12).Value = 12300
$ws.Cells.Item(137, 14).Value = -22500

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 300.45
$ws.Cells.Item(107, 9).Value = 295.1875
$ws.Cells.Item(107, 10).Value = 321.5
$ws.Cells.Item(107, 11).Value = 295.1875
$ws.Cells.Item(107, 12).Value = 321.5
$ws.Cells.Item(107, 13).Value = 1624.8125
$ws.Cells.Item(107, 14).Value = -4161.5
$ws.Cells.Item(113, 8).Value = 6766.0527
$ws.Cells.Item(113, 9).Value = 927.2222
$ws.Cells.Item(113, 10).Value = 12021
$ws.Cells.Item(113, 11).Value = 927.2222
$ws.Cells.Item(113, 12).Value = 12021
$ws.Cells.Item(113, 13).Value = 1242.7778
$ws.Cells.Item(113, 14).Value = -16361

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1606.1875
$ws.Cells.Item(7, 9).Value = 1366.6666
$ws.Cells.Item(7, 10).Value = 2324.75
$ws.Cells.Item(7, 11).Value = 1366.6666
$ws.Cells.Item(7, 12).Value = 2324.75
$ws.Cells.Item(7, 13).Value = -1254.6666
$ws.Cells.Item(7, 14).Value = -2548.75
$ws.Cells.Item(40, 8).Value = 4357.7915
$ws.Cells.Item(40, 9).Value = 3704.5557
$ws.Cells.Item(40, 10).Value = 6317.5
$ws.Cells.Item(40, 11).Value = 3704.5557
$ws.Cells.Item(40, 12).Value = 6317.5
$ws.Cells.Item(40, 13).Value = -3568.5557
$ws.Cells.Item(40, 14).Value = -6589.5
$ws.Cells.Item(110, 8).Value = 27661
$ws.Cells.Item(110, 10).Value = 27661
$ws.Cells.Item(110, 12).Value = 27661
$ws.Cells.Item(110, 14).Value = -35841
$ws.Cells.Item(126, 8).Value = 1606.1875
$ws.Cells.Item(126, 9).Value = 1366.6666
$ws.Cells.Item(126, 10).Value = 2324.75
$ws.Cells.Item(126, 11).Value = 4099.9998
$ws.Cells.Item(126, 12).Value = 6974.25
$ws.Cells.Item(126, 13).Value = -1629.9998
$ws.Cells.Item(126, 14).Value = -11914.25
$ws.Cells.Item(136, 8).Value = 1379.4667
$ws.Cells.Item(136, 9).Value = 1137.1177
$ws.Cells.Item(136, 10).Value = 2128.5454
$ws.Cells.Item(136, 11).Value = 3411.3531
$ws.Cells.Item(136, 12).Value = 6385.6362
$ws.Cells.Item(136, 13).Value = -861.3531000000003
$ws.Cells.Item(136, 14).Value = -11485.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 16685685
$ws.Cells.Item(62, 9).Value = 25026750
$ws.Cells.Item(62, 10).Value = 3555
$ws.Cells.Item(62, 11).Value = 25026750
$ws.Cells.Item(62, 12).Value = 3555
$ws.Cells.Item(62, 13).Value = -25026126
$ws.Cells.Item(62, 14).Value = -4803
$ws.Cells.Item(65, 8).Value = 16685685
$ws.Cells.Item(65, 9).Value = 25026750
$ws.Cells.Item(65, 10).Value = 3555
$ws.Cells.Item(65, 11).Value = 125133750
$ws.Cells.Item(65, 12).Value = 17775
$ws.Cells.Item(65, 13).Value = -125130630
$ws.Cells.Item(65, 14).Value = -24015
$ws.Cells.Item(113, 8).Value = 428.33334
$ws.Cells.Item(113, 9).Value = 317.5
$ws.Cells.Item(113, 10).Value = 650
$ws.Cells.Item(113, 11).Value = 952.5
$ws.Cells.Item(113, 12).Value = 1950
$ws.Cells.Item(113, 13).Value = 1217.5
$ws.Cells.Item(113, 14).Value = -6290

